$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the whole "State" design-pattern paragraph.
#    (commit message: "removed state")
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("A design pattern implemented in the final prototype")) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Merge the run-break (and drop the lastRenderedPageBreak marker) in the
#    middle of the Singleton paragraph so the sentence is one contiguous run
#    again. (commit message: "added singleton")
# ---------------------------------------------------------------------------
$mergeText = " in our system. Since our application is only running on one operating system, there is only instances of system objects. "
$mergeText += "Therefore, all of our methods access the system objects during runtime and do not need to account for third party systems. This allowed us to focus on implementation exclusively for Android OS. Without having to worry about multiple platforms, the group was able to write more robust code for Android. There wasn" + [char]0x2019 + "t a need for base functionality across platforms, so the system" + [char]0x2019 + "s code remains organized and clear."
$d.Content.Find.Execute($mergeText, $true, $false, $false, $false, $false, $true, 1, $false, $mergeText, 2)

# ---------------------------------------------------------------------------
# 3) Re-add a lastRenderedPageBreak marker right before the "Comparison of
#    Part 2 Class Diagram and Final Class Diagram" heading, since that is
#    now where pagination breaks after the deletion above.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Comparison of Part 2 Class Diagram and Final Class Diagram", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Collapse(1)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0C552810" w14:textId="35FC972D" w:rsidR="002D0074" w:rsidRDefault="002D0074"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Comparison of Part 2 Class Diagram and Final Class Diagram</w:t></w:r></w:p>'
    $r.InsertXML($xml)
}

Write-Output "edit complete"
